# Insert a new price-record row at row 425 (weekly "Acelga" report at
# Macroferia Regional de Talca), pushing all subsequent rows down by one.
# The previously-last row (499) now also occupies row 500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(425).Insert()

$ws.Range("A425").Value = 5
$ws.Range("B425").Value = "Macroferia Regional de Talca"
$ws.Range("C425").Value = "Maule"
$ws.Range("D425").Value = 45180
$ws.Range("E425").Value = 7
$ws.Range("F425").Value = 100112009
$ws.Range("G425").Value = "Acelga"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 500
$ws.Range("K425").Value = 2000
$ws.Range("L425").Value = 2000
$ws.Range("M425").Value = 2000
$ws.Range("N425").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O425").Value = "Provincia de Curicó"
$ws.Range("P425").Value = 500
$ws.Range("Q425").Value = 4
$ws.Range("R425").Value = "Hortaliza"
